$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.955.83'
$ws.Range('E2').Value = '  +3.42%  '

$ws.Range('D3').Value = '2.241.62'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '258.14'
$ws.Range('E5').Value = '  +2.28%  '

$ws.Range('D6').Value = '79.72'
$ws.Range('E6').Value = '  +8.06%  '

$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +2.09%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.600'
$ws.Range('E9').Value = '  +2.58%  '

$ws.Range('D10').Value = '43.14'
$ws.Range('E10').Value = '  +8.21%  '

$ws.Range('D11').Value = '0.0926'
$ws.Range('E11').Value = '  +1.13%  '

$ws.Range('D12').Value = '7.07'
$ws.Range('E12').Value = '  +4.02%  '

$ws.Range('D13').Value = '0.104'
$ws.Range('E13').Value = '  +2.31%  '

$ws.Range('D14').Value = '2.573.92'
$ws.Range('E14').Value = '  +1.56%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.70'
$ws.Range('E15').Value = '  +2.71%  '

$ws.Range('D16').Value = '2.227.31'
$ws.Range('E16').Value = '  +2.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.790'
$ws.Range('E17').Value = '  +1.97%  '

$ws.Range('D18').Value = '43.859.97'
$ws.Range('E18').Value = '  +3.30%  '

$ws.Range('E19').Value = '  +1.84%  '

$ws.Range('D20').Value = '71.37'
$ws.Range('E20').Value = '  +0.48%  '

$ws.Range('D21').Value = '6.03'
$ws.Range('E21').Value = '  +1.97%  '

$ws.Range('E22').Value = '  +6.47%  '

$ws.Range('D23').Value = '233.04'
$ws.Range('E23').Value = '  +2.12%  '

$ws.Range('E24').Value = '  -3.10%  '

$ws.Range('E25').Value = '  +0.13%  '

$ws.Range('D26').Value = '10.81'
$ws.Range('E26').Value = '  +1.82%  '

$ws.Range('D27').Value = '40.68'
$ws.Range('E27').Value = '  +9.30%  '

$ws.Range('E28').Value = '  +0.32%  '

$ws.Range('E29').Value = '  +0.84%  '

$ws.Range('E30').Value = '  -0.67%  '

$ws.Range('D31').Value = '172.54'
$ws.Range('E31').Value = '  +2.19%  '

$ws.Range('D32').Value = '0.0889'
$ws.Range('E32').Value = '  +11.00%  '

$ws.Range('D33').Value = '20.58'
$ws.Range('E33').Value = '  +2.72%  '

$ws.Range('D34').Value = '5.29'
$ws.Range('E34').Value = '  +2.57%  '

$ws.Range('E35').Value = '  +2.11%  '

$ws.Range('E36').Value = '  +5.24%  '

$ws.Range('D37').Value = '0.0366'
$ws.Range('E37').Value = '  +11.02%  '

$ws.Range('E38').Value = '  +3.64%  '

$ws.Range('D39').Value = '12.97'
$ws.Range('E39').Value = '  +7.59%  '

$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  +23.67%  '

$ws.Range('D41').Value = '2.14'
$ws.Range('E41').Value = '  +3.44%  '

$ws.Range('D42').Value = '5.52'
$ws.Range('E42').Value = '  +5.29%  '

$ws.Range('D43').Value = '62.84'
$ws.Range('E43').Value = '  +6.51%  '

$ws.Range('D44').Value = '0.202'
$ws.Range('E44').Value = '  +2.40%  '

$ws.Range('D45').Value = '103.94'
$ws.Range('E45').Value = '  +1.44%  '

$ws.Range('D46').Value = '8.51'
$ws.Range('E46').Value = '  +2.20%  '

$ws.Range('D47').Value = '0.0985'
$ws.Range('E47').Value = '  +0.82%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '1.12'
$ws.Range('E48').Value = '  +2.25%  '

$ws.Range('B49').Value = 'WOONetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D49').Value = '0.447'
$ws.Range('E49').Value = '  -7.08%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').Value = '  +26.76%  '

$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '1.15'
$ws.Range('E51').Value = '  +1.70%  '
